# Add the CanESM5 (CCCma) atmosphere-physics citation entries to the
# "Citations" sheet, tidy up that sheet's banded-row formatting so the
# light/dark zebra-striping boundaries still line up, drop two now-unused
# trailing rows, and leave the workbook with the Citations tab active
# (mirrors the author re-opening the template, filling in the citations,
# and saving).

$wb = $excel.ActiveWorkbook

$wsFrontis   = $wb.Worksheets.Item("Frontis")
$wsExample   = $wb.Worksheets.Item("Example")
$wsCitations = $wb.Worksheets.Item("Citations")

# ---------------------------------------------------------------------
# 1. Populate the citation table (columns A = Identifier/Mnemonic,
#    B = DOI) on the Citations sheet, rows 3-15.
# ---------------------------------------------------------------------

$wsCitations.Cells.Item(3, 1).Value = "Scinocca_et_al_2008"
$wsCitations.Cells.Item(3, 2).Value = "10.5194/acp-8-7055-2008"

$wsCitations.Cells.Item(4, 1).Value = "von_Salzen_et_al_2013"
$wsCitations.Cells.Item(4, 2).Value = "10.1080/07055900.2012.755610"

$wsCitations.Cells.Item(5, 2).Value = "10.1080/07055900.1995.9649539`n"
$wsCitations.Cells.Item(5, 1).Value = "Zhang_McFarlane_1995"

$wsCitations.Cells.Item(6, 1).Value = "Scinocca_McFarlane_2004"
$wsCitations.Cells.Item(6, 2).Value = "10.1175/1520-0469(2004)061<1993:TVOMTP>2.0.CO;2"

$wsCitations.Cells.Item(7, 1).Value = "von_Salzen_McFarlane_2002"

$wsCitations.Cells.Item(8, 2).Value = "10.1007/BF00207939"

$wsCitations.Cells.Item(7, 2).Value = "10.1175/1520-0469(2002)059<1405:POTBEO>2.0.CO;2"

$wsCitations.Cells.Item(8, 1).Value = "Lohmann_Roeckner_1996"

$wsCitations.Cells.Item(9, 2).Value = "10.1002/qj.49712656914"

$wsCitations.Cells.Item(9, 1).Value = "Hogan_Illingworth_2000"

$wsCitations.Cells.Item(10, 1).Value = "Scinocca_McFarlane_2000"
$wsCitations.Cells.Item(10, 2).Value = "10.1002/qj.49712656802"

$wsCitations.Cells.Item(11, 2).Value = "10.1175/1520-0469(1987)044<1775:TEOOEG>2.0.CO;2"

$wsCitations.Cells.Item(11, 1).Value = "McFarlane_1987"

$wsCitations.Cells.Item(12, 1).Value = "Arora_Boer_1999"
$wsCitations.Cells.Item(12, 2).Value = "10.1029/1999JD900905`n"

$wsCitations.Cells.Item(13, 1).Value = "Lohmann_et_al_1999"
$wsCitations.Cells.Item(13, 2).Value = "10.1029/1999JD900343"

$wsCitations.Cells.Item(14, 1).Value = "Croft_et_al_2005"
$wsCitations.Cells.Item(14, 2).Value = "10.5194/acp-5-1931-2005"

$wsCitations.Cells.Item(15, 1).Value = "Lana_et_al_2010"
$wsCitations.Cells.Item(15, 2).Value = "10.1029/2010GB003850"

# ---------------------------------------------------------------------
# 2. Remove two now-superfluous blank rows from the banded filler area
#    (this also re-aligns the light/dark striping boundary and shrinks
#    the sheet's used range).
# ---------------------------------------------------------------------

$wsCitations.Rows.Item(20).Delete()
$wsCitations.Rows.Item(20).Delete()

# ---------------------------------------------------------------------
# 3. Minor row-height retouch left behind by the resave (header band
#    and the uniform body rows on Frontis/Example/Citations).
# ---------------------------------------------------------------------

$wsFrontis.Rows.Item(1).RowHeight = 34.35
$wsFrontis.Rows.Item(6).RowHeight = 408.95
$wsFrontis.Rows.Item(7).RowHeight = 408.95

$wsExample.Rows.Item(1).RowHeight = 30.95
$wsExample.Rows.Item(2).RowHeight = 20.85
for ($r = 4; $r -le 56; $r++) {
    $wsExample.Rows.Item($r).RowHeight = 20.1
}

$wsCitations.Rows.Item(1).RowHeight = 30.95
$wsCitations.Rows.Item(2).RowHeight = 20.85
for ($r = 4; $r -le 50; $r++) {
    $wsCitations.Rows.Item($r).RowHeight = 20.1
}

# ---------------------------------------------------------------------
# 4. Leave the view the way the author left it: scrolled/selected on
#    Frontis!B5 first, then finally parked on Citations!B15 (so the
#    Citations tab is the active tab on save).
# ---------------------------------------------------------------------

$wsFrontis.Activate()
$wsFrontis.Range("B5").Select()

$wsCitations.Activate()
$wsCitations.Range("B15").Select()

Write-Host "done"
